$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B176").Formula = "=0.7019"
